# The site generator stopped emitting the trailing "Ver no Jupiter /
# Salvar em pdf / Salvar em docx" line and the "(c) 2020 ... Creative
# Commons" footer line on this page. Remove both paragraphs, together
# with the blank paragraph that used to separate them from the
# "Requisitos" section above, while leaving everything else (including
# the blank paragraph + page-break paragraph that follow) untouched.

$d = $word.ActiveDocument

$jupIdx = 0
$copyrightIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $jupIdx = $i
    }
    if ($t -like "*Creative Commons*") {
        $copyrightIdx = $i
    }
}

if ($jupIdx -gt 0 -and $copyrightIdx -ge $jupIdx) {
    # Also sweep up the blank paragraph immediately before the
    # "Ver no Jupiter..." paragraph.
    $startPara = $d.Paragraphs.Item($jupIdx - 1)
    $endPara   = $d.Paragraphs.Item($copyrightIdx)

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
